# Updated symbol list (Price / Volume(1h) columns) to reflect the latest
# crypto snapshot, as produced by the scheduled GitHub Actions scraper.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Collect the set of cells we need to update and force them to Text format
# so Excel does not auto-convert the numeric-looking strings into numbers.
$targetCells = @("D2","E2","D3","E3","D4","E4","D5","D6","E6","D7","E7","D8","E8","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","E17","D18","E18","D19","E19","E20","D21","E21","E23","D24","E24","D25","E25","D26","E26","D39","E39","D40","E40","D41","E41","D42","E42","E43","D44","E44","D45","E45","D46","E47","D48","E48","D49","E49","D50","E50","D51","E51")
foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "324.92"
$ws.Range("E2").Value = "-1.32%"
$ws.Range("D3").Value = "39.67"
$ws.Range("E3").Value = "-0.80%"
$ws.Range("D4").Value = "5.627"
$ws.Range("E4").Value = "6.03%"
$ws.Range("D5").Value = "0.08017"
$ws.Range("D6").Value = "2.015"
$ws.Range("E6").Value = "4.70%"
$ws.Range("D7").Value = "4.486"
$ws.Range("E7").Value = "-0.74%"
$ws.Range("D8").Value = "8.619"
$ws.Range("E8").Value = "-0.37%"
$ws.Range("E9").Value = "-1.58%"
$ws.Range("D10").Value = "0.9225"
$ws.Range("E10").Value = "-2.21%"
$ws.Range("D11").Value = "0.1242"
$ws.Range("E11").Value = "-8.41%"
$ws.Range("D12").Value = "0.1956"
$ws.Range("E12").Value = "-0.81%"
$ws.Range("D13").Value = "8.725"
$ws.Range("E13").Value = "21.01%"
$ws.Range("D14").Value = "0.09109"
$ws.Range("E14").Value = "-2.32%"
$ws.Range("D15").Value = "0.03571"
$ws.Range("E15").Value = "0.51%"
$ws.Range("D16").Value = "0.1048"
$ws.Range("E16").Value = "9.28%"
$ws.Range("E17").Value = "-2.39%"
$ws.Range("D18").Value = "0.006244"
$ws.Range("E18").Value = "2.96%"
$ws.Range("D19").Value = "3.351"
$ws.Range("E19").Value = "-0.38%"
$ws.Range("E20").Value = "-0.92%"
$ws.Range("D21").Value = "0.1372"
$ws.Range("E21").Value = "3.07%"
$ws.Range("E23").Value = "-1.36%"
$ws.Range("D24").Value = "0.001262"
$ws.Range("E24").Value = "3.27%"
$ws.Range("D25").Value = "0.004604"
$ws.Range("E25").Value = "6.58%"
$ws.Range("D26").Value = "0.0001231"
$ws.Range("E26").Value = "2.45%"
$ws.Range("D39").Value = "0.02488"
$ws.Range("E39").Value = "0.08%"
$ws.Range("D40").Value = "0.05320"
$ws.Range("E40").Value = "2.11%"
$ws.Range("D41").Value = "0.007477"
$ws.Range("E41").Value = "-3.49%"
$ws.Range("D42").Value = "0.009905"
$ws.Range("E42").Value = "7.86%"
$ws.Range("E43").Value = "-1.87%"
$ws.Range("D44").Value = "0.002117"
$ws.Range("E44").Value = "-2.09%"
$ws.Range("D45").Value = "0.01114"
$ws.Range("E45").Value = "13.13%"
$ws.Range("D46").Value = "0.00006693"
$ws.Range("E47").Value = "-0.04%"
$ws.Range("D48").Value = "0.002976"
$ws.Range("E48").Value = "-11.01%"
$ws.Range("D49").Value = "0.002282"
$ws.Range("E49").Value = "-5.03%"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").Value = "-0.04%"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").Value = "-0.04%"

# Restore default (Normal) style so no stray number-format style is left behind
foreach ($addr in $targetCells) {
    $ws.Range($addr).Style = "Normal"
}
